$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume/ordering refresh (GitHub Actions daily scrape update).
# D (Price) and E (Volume 1h) columns hold text that often *looks* numeric or
# percent-like (e.g. "244.95", "1,594.13%"). Excel COM auto-converts such strings
# to real numbers/percentages on assignment, which would corrupt the original
# text-cell semantics of this sheet. To keep them as plain text we temporarily
# force the Text number format, assign the literal string, then reset the style
# back to Normal so no stray formatting is left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "244.95"
Set-TextValue $ws.Range("E2") "-0.68%"
Set-TextValue $ws.Range("D3") "27.11"
Set-TextValue $ws.Range("E3") "3.17%"
Set-TextValue $ws.Range("D4") "5.088"
Set-TextValue $ws.Range("E4") "0.02%"
Set-TextValue $ws.Range("D5") "0.05701"
Set-TextValue $ws.Range("E5") "1.78%"
Set-TextValue $ws.Range("D6") "6.506"
Set-TextValue $ws.Range("E6") "0.40%"
Set-TextValue $ws.Range("D7") "0.8203"
Set-TextValue $ws.Range("E7") "0.82%"
Set-TextValue $ws.Range("D8") "0.8631"
Set-TextValue $ws.Range("E8") "1.89%"
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D9") "0.01010"
Set-TextValue $ws.Range("E9") "1,594.13%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D10") "0.1332"
Set-TextValue $ws.Range("E10") "-0.74%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D11") "0.06907"
Set-TextValue $ws.Range("E11") "-1.45%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D12") "0.02828"
Set-TextValue $ws.Range("E12") "-1.03%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D13") "0.09372"
Set-TextValue $ws.Range("E13") "-0.18%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D14") "0.001518"
Set-TextValue $ws.Range("E14") "-0.11%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D15") "0.04084"
Set-TextValue $ws.Range("E15") "-12.25%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D16") "0.006040"
Set-TextValue $ws.Range("E16") "-1.93%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D17") "3.508"
Set-TextValue $ws.Range("E17") "-2.42%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D18") "3.010"
Set-TextValue $ws.Range("E18") "-0.02%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D19") "2.315"
Set-TextValue $ws.Range("E19") "12.61%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue $ws.Range("D20") "0.3177"
Set-TextValue $ws.Range("E20") "-0.93%"
$ws.Range("B21").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C21").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D21") "0.03173"
Set-TextValue $ws.Range("E21") "-0.96%"
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue $ws.Range("D22") "0.1274"
Set-TextValue $ws.Range("E22") "-1.79%"
$ws.Range("B23").Value = "MCDex"
$ws.Range("C23").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D23") "3.565"
Set-TextValue $ws.Range("E23") "-4.81%"
$ws.Range("B24").Value = "ZBToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue $ws.Range("D24") "0.1374"
Set-TextValue $ws.Range("E24") "1.75%"
Set-TextValue $ws.Range("D25") "0.001217"
Set-TextValue $ws.Range("E25") "-2.22%"
Set-TextValue $ws.Range("D26") "0.003972"
Set-TextValue $ws.Range("E26") "-13.32%"
Set-TextValue $ws.Range("D27") "0.00009899"
Set-TextValue $ws.Range("E27") "3.10%"
Set-TextValue $ws.Range("D28") "0.0001449"
Set-TextValue $ws.Range("E28") "-25.26%"
Set-TextValue $ws.Range("D40") "0.03721"
Set-TextValue $ws.Range("E40") "1.47%"
Set-TextValue $ws.Range("D41") "0.005709"
Set-TextValue $ws.Range("E41") "-7.78%"
Set-TextValue $ws.Range("D42") "0.1057"
Set-TextValue $ws.Range("E42") "0.12%"
Set-TextValue $ws.Range("D43") "0.002368"
Set-TextValue $ws.Range("E43") "-5.30%"
Set-TextValue $ws.Range("D44") "0.009372"
Set-TextValue $ws.Range("E44") "6.80%"
Set-TextValue $ws.Range("D45") "0.00005169"
Set-TextValue $ws.Range("E45") "-2.40%"
Set-TextValue $ws.Range("E46") "-0.03%"
Set-TextValue $ws.Range("D47") "0.1015"
Set-TextValue $ws.Range("E47") "-7.76%"
Set-TextValue $ws.Range("D48") "0.002537"
Set-TextValue $ws.Range("E48") "-4.77%"
Set-TextValue $ws.Range("D49") "0.00002100"
Set-TextValue $ws.Range("E49") "-0.03%"
Set-TextValue $ws.Range("D50") "0.0002000"
Set-TextValue $ws.Range("E50") "-0.03%"
